# Generate Report for Handoff
#
# Updates the localization-status report after a new handoff report run:
#   - The "Priority" column (E) for the rows whose handback validation
#     failed with a handoff-type mismatch now records the handoff type "ht"
#     (previously blank) on both the zh-cn and de-de sheets.
#   - The "Latest Handoff Datetime" column (H) for those same rows is
#     refreshed to the timestamp of the new handoff run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$rows = @(7, 9, 10, 11, 13, 14)

foreach ($row in $rows) {
    $wsZhCn.Range("E$row").Value = "ht"
    $wsZhCn.Range("H$row").Value = "2016-09-01 14:26:41"

    $wsDeDe.Range("E$row").Value = "ht"
    $wsDeDe.Range("H$row").Value = "2016-09-01 14:26:47"
}
